# Update the "D1"/"D2"/"D3" detail-angle labels (previously all "細節")
# and fix the "_Fex__" typo (trailing double underscore) to "_Fex_".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 13-15 already use the lighter "theme" font style; just change the text.
$ws.Range("C13").Value = "D1"
$ws.Range("C14").Value = "D2"
$ws.Range("C15").Value = "D3"

# Rows 16-18 pick up the same font formatting as C13:C15 (copy formats) before
# the text is switched, matching how these duplicate labels were normalized.
$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C18").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C16").Value = "D1"
$ws.Range("C17").Value = "D2"
$ws.Range("C18").Value = "D3"

$ws.Range("A33").Value = "_Fex_"

# Leave the selection where the author last left it before saving.
$ws.Range("H33").Select() | Out-Null
